$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.935.06'
$ws.Range('E2').Value = '  +0.69%  '

$ws.Range('D3').Value = '2.531.30'
$ws.Range('E3').Value = '  +0.38%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.44'
$ws.Range('D5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.69'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.579'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.17%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.98%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.02%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0811'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.05%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.60'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.97%  '

$ws.Range('E13').Value = '  -0.49%  '

$ws.Range('D14').Value = '2.920.92'
$ws.Range('E14').Value = '  +0.51%  '

$ws.Range('D15').Value = '2.587.39'
$ws.Range('E15').Value = '  +3.58%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.86%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.850'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.46%  '

$ws.Range('D18').Value = '43.016.42'
$ws.Range('E18').Value = '  +0.81%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.03'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.90%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.12%  '

$ws.Range('D21').Value = '0.0₃0967'
$ws.Range('E21').Value = '  -0.63%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.41%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '251.79'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.44%  '

$ws.Range('E25').Value = '  -0.71%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.25%  '

$ws.Range('E27').Value = '  +0.00%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.42'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.36%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.85'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.69%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.01%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.03'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.32'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.31%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.72%  '

$ws.Range('E34').Value = '  +0.10%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0791'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.12%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.80'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.14%  '

$ws.Range('E37').Value = '  -0.52%  '

$ws.Range('E38').Value = '  -2.73%  '

$ws.Range('E39').Value = '  -0.37%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.70%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +14.45%  '

$ws.Range('E44').Value = '  +0.37%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.29'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.97%  '

$ws.Range('D46').Value = '2.022.60'
$ws.Range('E46').Value = '  -0.43%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.91'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.56%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.80'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.79%  '

$ws.Range('D49').Value = '2.776.04'
$ws.Range('E49').Value = '  +0.40%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.68'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.77%  '

# Row 42/43: RenderToken and VeChain swap places, with updated prices/volumes
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0305'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.71%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.81'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.54%  '
